$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.117734
$ws.Range("H2").Value = 0.353202
$ws.Range("I2").Value = 0.6076081328197709
$ws.Range("J2").Value = 0.6076081328197709
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8155003333333334
$ws.Range("N2").Value = 2.446501
$ws.Range("O2").Value = 0.1910612426590028
$ws.Range("P2").Value = 0.1910612426590029
$ws.Range("Q2").Value = 0.09601211624466667
$ws.Range("R2").Value = 0.864109046202
$ws.Range("S2").Value = 0.1160903649062619
$ws.Range("T2").Value = 0.1160903649062619

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.117734
$ws.Range("H3").Value = 0.353202
$ws.Range("I3").Value = 0.6076081328197709
$ws.Range("J3").Value = 0.6076081328197709
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.333134333333334
$ws.Range("N3").Value = 9.999403000000001
$ws.Range("O3").Value = 0.7809105179307759
$ws.Range("P3").Value = 0.780910517930776
$ws.Range("Q3").Value = 0.3924232376006667
$ws.Range("R3").Value = 3.531809138406
$ws.Range("S3").Value = 0.474487581699239
$ws.Range("T3").Value = 0.474487581699239

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.117734
$ws.Range("H4").Value = 0.353202
$ws.Range("I4").Value = 0.6076081328197709
$ws.Range("J4").Value = 0.6076081328197709
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.119632
$ws.Range("N4").Value = 0.358896
$ws.Range("O4").Value = 0.02802823941022116
$ws.Range("P4").Value = 0.02802823941022117
$ws.Range("Q4").Value = 0.014084753888
$ws.Range("R4").Value = 0.126762784992
$ws.Range("S4").Value = 0.01703018621427
$ws.Range("T4").Value = 0.01703018621427

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf15"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.029478
$ws.Range("H5").Value = 0.088434
$ws.Range("I5").Value = 0.1521316912638762
$ws.Range("J5").Value = 0.1521316912638762
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8155003333333334
$ws.Range("N5").Value = 2.446501
$ws.Range("O5").Value = 0.1910612426590028
$ws.Range("P5").Value = 0.1910612426590029
$ws.Range("Q5").Value = 0.024039318826
$ws.Range("R5").Value = 0.216353869434
$ws.Range("S5").Value = 0.02906646998069196
$ws.Range("T5").Value = 0.02906646998069197

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf15"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.029478
$ws.Range("H6").Value = 0.088434
$ws.Range("I6").Value = 0.1521316912638762
$ws.Range("J6").Value = 0.1521316912638762
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.333134333333334
$ws.Range("N6").Value = 9.999403000000001
$ws.Range("O6").Value = 0.7809105179307759
$ws.Range("P6").Value = 0.780910517930776
$ws.Range("Q6").Value = 0.09825413387800001
$ws.Range("R6").Value = 0.8842872049020001
$ws.Range("S6").Value = 0.1188012378185585
$ws.Range("T6").Value = 0.1188012378185585

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf15"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.029478
$ws.Range("H7").Value = 0.088434
$ws.Range("I7").Value = 0.1521316912638762
$ws.Range("J7").Value = 0.1521316912638762
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.119632
$ws.Range("N7").Value = 0.358896
$ws.Range("O7").Value = 0.02802823941022116
$ws.Range("P7").Value = 0.02802823941022117
$ws.Range("Q7").Value = 0.003526512096
$ws.Range("R7").Value = 0.031738608864
$ws.Range("S7").Value = 0.004263983464625774
$ws.Range("T7").Value = 0.004263983464625775

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf15"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.04655433333333334
$ws.Range("H8").Value = 0.139663
$ws.Range("I8").Value = 0.2402601759163528
$ws.Range("J8").Value = 0.2402601759163529
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8155003333333334
$ws.Range("N8").Value = 2.446501
$ws.Range("O8").Value = 0.1910612426590028
$ws.Range("P8").Value = 0.1910612426590029
$ws.Range("Q8").Value = 0.03796507435144445
$ws.Range("R8").Value = 0.341685669163
$ws.Range("S8").Value = 0.045904407772049
$ws.Range("T8").Value = 0.04590440777204902

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf15"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.04655433333333334
$ws.Range("H9").Value = 0.139663
$ws.Range("I9").Value = 0.2402601759163528
$ws.Range("J9").Value = 0.2402601759163529
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.333134333333334
$ws.Range("N9").Value = 9.999403000000001
$ws.Range("O9").Value = 0.7809105179307759
$ws.Range("P9").Value = 0.780910517930776
$ws.Range("Q9").Value = 0.1551718467987778
$ws.Range("R9").Value = 1.396546621189
$ws.Range("S9").Value = 0.1876216984129784
$ws.Range("T9").Value = 0.1876216984129785

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf15"
$ws.Range("C10").Value = "Fgfr2"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.04655433333333334
$ws.Range("H10").Value = 0.139663
$ws.Range("I10").Value = 0.2402601759163528
$ws.Range("J10").Value = 0.2402601759163529
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.119632
$ws.Range("N10").Value = 0.358896
$ws.Range("O10").Value = 0.02802823941022116
$ws.Range("P10").Value = 0.02802823941022117
$ws.Range("Q10").Value = 0.005569388005333334
$ws.Range("R10").Value = 0.050124492048
$ws.Range("S10").Value = 0.006734069731325391
$ws.Range("T10").Value = 0.006734069731325392

